# Update odds figures per the 2025-03-03 FlashScore refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q3").Value = 1.73
$ws.Range("R3").Value = 2.1
$ws.Range("M4").Value = 1.11
$ws.Range("N4").Value = 6.5
$ws.Range("AR6").Value = 3.95
$ws.Range("U7").Value = 1.78
$ws.Range("V7").Value = 2.03
$ws.Range("G9").Value = 10
$ws.Range("H9").Value = 5
$ws.Range("I9").Value = 1.26
$ws.Range("J9").Value = 8.25
$ws.Range("K9").Value = 2.42
$ws.Range("L9").Value = 1.7
$ws.Range("O9").Value = 1.21
$ws.Range("P9").Value = 3.5
$ws.Range("Q9").Value = 1.65
$ws.Range("R9").Value = 2
$ws.Range("S9").Value = 2.52
$ws.Range("T9").Value = 1.4
$ws.Range("W9").Value = 2.15
$ws.Range("X9").Value = 1.55
$ws.Range("Y9").Value = 24
$ws.Range("Z9").Value = 75
$ws.Range("AA9").Value = 32
$ws.Range("AB9").Value = 350
$ws.Range("AD9").Value = 120
$ws.Range("AE9").Value = 12
$ws.Range("AF9").Value = 10.5
$ws.Range("AG9").Value = 27
$ws.Range("AH9").Value = 150
$ws.Range("AI9").Value = 6.5
$ws.Range("AJ9").Value = 5.7
$ws.Range("AK9").Value = 9
$ws.Range("AL9").Value = 7.2
$ws.Range("AM9").Value = 11.5
$ws.Range("AN9").Value = 35
$ws.Range("M11").Value = 1.05
$ws.Range("N11").Value = 11
$ws.Range("Q11").Value = 2.05
$ws.Range("R11").Value = 1.75
$ws.Range("S11").Value = 3.5
$ws.Range("T11").Value = 1.29
$ws.Range("M12").Value = 1.01
$ws.Range("O12").Value = 1.08
$ws.Range("S12").Value = 1.87
$ws.Range("T12").Value = 1.77
$ws.Range("M13").Value = 1.07
$ws.Range("O13").Value = 1.47
$ws.Range("T13").Value = 1.13
$ws.Range("AP13").Value = 1.92
$ws.Range("AQ13").Value = 1.82
$ws.Range("M14").Value = 1.05
$ws.Range("N14").Value = 9
$ws.Range("O14").Value = 1.33
$ws.Range("Q14").Value = 2.2
$ws.Range("R14").Value = 1.65
$ws.Range("T14").Value = 1.19
$ws.Range("M15").Value = 1.02
$ws.Range("O15").Value = 1.17
$ws.Range("S15").Value = 2.62
$ws.Range("T15").Value = 1.41
$ws.Range("M16").Value = 1.05
$ws.Range("O16").Value = 1.33
$ws.Range("T16").Value = 1.19
$ws.Range("AP20").Value = 2.1
$ws.Range("AQ20").Value = 1.78
$ws.Range("G21").Value = 2.55
$ws.Range("I21").Value = 2.63
$ws.Range("L21").Value = 3.2
$ws.Range("G22").Value = 2.15
$ws.Range("I22").Value = 3.1
$ws.Range("J22").Value = 2.72
$ws.Range("K22").Value = 2.1
$ws.Range("L22").Value = 3.6
$ws.Range("O22").Value = 1.23
$ws.Range("P22").Value = 3.35
$ws.Range("S22").Value = 2.62
$ws.Range("T22").Value = 1.38
$ws.Range("X22").Value = 2.12
$ws.Range("Z22").Value = 11.75
$ws.Range("AA22").Value = 8.5
$ws.Range("AH22").Value = 45
$ws.Range("AI22").Value = 11
$ws.Range("AJ22").Value = 18
$ws.Range("AK22").Value = 10.75
$ws.Range("AN22").Value = 28
$ws.Range("Q25").Value = 1.88
$ws.Range("R25").Value = 1.98
$ws.Range("S25").Value = 3.25
$ws.Range("T25").Value = 1.33
$ws.Range("G26").Value = 2.55
$ws.Range("I26").Value = 2.8
$ws.Range("L26").Value = 3.6
$ws.Range("Y26").Value = 7
$ws.Range("Z26").Value = 11
$ws.Range("AJ26").Value = 13
$ws.Range("AM26").Value = 26
$ws.Range("AN26").Value = 41
$ws.Range("G27").Value = 2.05
$ws.Range("H27").Value = 3.1
$ws.Range("I27").Value = 3.4
$ws.Range("J27").Value = 2.65
$ws.Range("K27").Value = 2.14
$ws.Range("L27").Value = 4.1
$ws.Range("O27").Value = 1.32
$ws.Range("P27").Value = 3.1
$ws.Range("Q27").Value = 1.99
$ws.Range("R27").Value = 1.76
$ws.Range("S27").Value = 3.4
$ws.Range("T27").Value = 1.29
$ws.Range("U27").Value = 1.42
$ws.Range("V27").Value = 2.65
$ws.Range("W27").Value = 1.79
$ws.Range("X27").Value = 1.92
$ws.Range("Y27").Value = 6.8
$ws.Range("Z27").Value = 8.2
$ws.Range("AA27").Value = 7.6
$ws.Range("AB27").Value = 16
$ws.Range("AC27").Value = 15
$ws.Range("AD27").Value = 27
$ws.Range("AE27").Value = 9.4
$ws.Range("AF27").Value = 5.8
$ws.Range("AG27").Value = 12
$ws.Range("AH27").Value = 70
$ws.Range("AI27").Value = 9.4
$ws.Range("AJ27").Value = 15
$ws.Range("AK27").Value = 10
$ws.Range("AL27").Value = 40
$ws.Range("AM27").Value = 28
$ws.Range("AN27").Value = 35
$ws.Range("AO27").Value = 101
$ws.Range("K28").Value = 2.25
$ws.Range("O28").Value = 1.22
$ws.Range("P28").Value = 4
$ws.Range("Q28").Value = 1.75
$ws.Range("R28").Value = 2.05
$ws.Range("S28").Value = 2.75
$ws.Range("T28").Value = 1.4
$ws.Range("Y28").Value = 9
$ws.Range("AK28").Value = 12
$ws.Range("O30").Value = 1.18
$ws.Range("P30").Value = 4.5
$ws.Range("Q30").Value = 1.6
$ws.Range("R30").Value = 2.3
$ws.Range("S30").Value = 2.5
$ws.Range("T30").Value = 1.5
$ws.Range("AR30").Value = 2
$ws.Range("AS30").Value = 1.85
$ws.Range("G43").Value = 2.35
$ws.Range("I43").Value = 2.9
$ws.Range("J43").Value = 3
$ws.Range("L43").Value = 3.5
$ws.Range("AM43").Value = 23
